$d = $word.ActiveDocument

# Replace "konstelasi sepatu bot" with "Rasi sepatu bot" everywhere it occurs.
# (Other standalone occurrences of "konstelasi" / "konstelasimu" elsewhere in
# the document are unrelated and must be left untouched.)
$d.Content.Find.Execute("konstelasi sepatu bot", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Rasi sepatu bot", 2)
